# Insert a new weekly price record for Cilantro (Feria Lagunitas de Puerto Montt)
# as row 410, shifting existing rows 410-437 down to 411-438.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 410 (pushes old rows 410..437 -> 411..438,
# carries formatting down the same way Excel's native row-insert does).
$ws.Rows.Item(410).Insert()

# Populate the newly inserted row 410 with the new record.
$ws.Cells.Item(410, 1).Value  = 4
$ws.Cells.Item(410, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(410, 3).Value  = "Los Lagos"
$ws.Cells.Item(410, 4).Value  = 45021
$ws.Cells.Item(410, 5).Value  = 10
$ws.Cells.Item(410, 6).Value  = 100112040
$ws.Cells.Item(410, 7).Value  = "Cilantro"
$ws.Cells.Item(410, 8).Value  = "Sin especificar"
$ws.Cells.Item(410, 9).Value  = "Primera"
$ws.Cells.Item(410, 10).Value = 20
$ws.Cells.Item(410, 11).Value = 8000
$ws.Cells.Item(410, 12).Value = 8000
$ws.Cells.Item(410, 13).Value = 8000
$ws.Cells.Item(410, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(410, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(410, 16).Value = 4000
$ws.Cells.Item(410, 17).Value = 2
$ws.Cells.Item(410, 18).Value = "Hortaliza"
